$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '29.248.83'
$ws.Cells.Item(2, 5).Value = '  +0.39%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.860.49'
$ws.Cells.Item(3, 5).Value = '  +0.62%  '

$ws.Cells.Item(4, 5).Value = '  +0.07%  '

$ws.Cells.Item(5, 5).Value = '  -0.84%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '237.44'
$ws.Cells.Item(6, 5).Value = '  -0.30%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.002'
$ws.Cells.Item(7, 5).Value = '  +0.11%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.08149'

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.3025'
$ws.Cells.Item(9, 5).Value = '  -0.40%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '23.19'
$ws.Cells.Item(10, 5).Value = '  -0.85%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.08159'
$ws.Cells.Item(11, 5).Value = '  +0.35%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.854.42'
$ws.Cells.Item(12, 5).Value = '  +0.15%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.154'
$ws.Cells.Item(13, 5).Value = '  -1.08%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.7058'
$ws.Cells.Item(14, 5).Value = '  -2.68%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '89.02'
$ws.Cells.Item(15, 5).Value = '  +0.24%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '29.257.66'
$ws.Cells.Item(16, 5).Value = '  +0.22%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '5.764'
$ws.Cells.Item(17, 5).Value = '  +0.26%  '

$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000007824'
$ws.Cells.Item(18, 5).Value = '  +2.31%  '

$ws.Cells.Item(19, 2).Value = 'Avalanche'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.31'
$ws.Cells.Item(19, 5).Value = '  +1.85%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '235.71'
$ws.Cells.Item(20, 5).Value = '  -1.14%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.001'
$ws.Cells.Item(21, 5).Value = '  +0.11%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '2.112.63'
$ws.Cells.Item(22, 5).Value = '  +0.84%  '

$ws.Cells.Item(23, 5).Value = '  +0.11%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '7.418'
$ws.Cells.Item(24, 5).Value = '  -1.73%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '161.61'
$ws.Cells.Item(25, 5).Value = '  +0.09%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '8.947'
$ws.Cells.Item(26, 5).Value = '  -0.40%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.1439'
$ws.Cells.Item(27, 5).Value = '  -1.23%  '

$ws.Cells.Item(28, 5).Value = '  +0.32%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.959'
$ws.Cells.Item(29, 5).Value = '  +0.23%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.432'
$ws.Cells.Item(30, 5).Value = '  +3.08%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.480'
$ws.Cells.Item(31, 5).Value = '  -1.02%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.393'
$ws.Cells.Item(32, 5).Value = '  -2.48%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.052'
$ws.Cells.Item(33, 5).Value = '  +1.99%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.05190'
$ws.Cells.Item(34, 5).Value = '  +0.75%  '

$ws.Cells.Item(35, 5).Value = '  -1.48%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.7066'
$ws.Cells.Item(36, 5).Value = '  +1.04%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.000'
$ws.Cells.Item(37, 5).Value = '  -3.71%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.677'
$ws.Cells.Item(38, 5).Value = '  +0.82%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01835'
$ws.Cells.Item(39, 5).Value = '  -2.00%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.727'
$ws.Cells.Item(40, 5).Value = '  +1.70%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.9219'
$ws.Cells.Item(41, 5).Value = '  -1.00%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.135.94'
$ws.Cells.Item(42, 5).Value = '  +4.87%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.4257'
$ws.Cells.Item(43, 5).Value = '  -0.60%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '5.892'
$ws.Cells.Item(44, 5).Value = '  -1.77%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '70.22'
$ws.Cells.Item(45, 5).Value = '  +0.72%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.001'
$ws.Cells.Item(46, 5).Value = '  +0.10%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '102.22'
$ws.Cells.Item(47, 5).Value = '  +0.12%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.764'
$ws.Cells.Item(48, 5).Value = '  +1.29%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.006.97'
$ws.Cells.Item(49, 5).Value = '  +0.70%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '9.148'
$ws.Cells.Item(50, 5).Value = '  +0.07%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '6.935'
$ws.Cells.Item(51, 5).Value = '  -1.38%  '
